$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.121.50"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.791.48"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "229.07"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "0.552"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "32.58"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "0.288"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").Value = "0.0713"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.049.65"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "11.09"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "1.792.06"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "0.626"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "34.065.94"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "68.60"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "245.57"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "160.67"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "7.11"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "1.26"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("D31").Value = "0.0516"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").Value = "3.52"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").Value = "1.398.11"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "0.662"
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "0.917"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "78.57"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "13.22"
$ws.Range("E44").Value = "  +10.15%  "
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("E46").Value = "  +10.27%  "
$ws.Range("D47").Value = "109.30"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "5.85"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").Value = "1.949.02"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.19%  "
